$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (10) of notification data, matching the existing table layout.
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Thông báo 1"
$ws.Cells.Item(10, 3).Value = "<p>123</p>"
$ws.Cells.Item(10, 4).Value = "Ban Đào Tạo"
$ws.Cells.Item(10, 5).Value = "19/09/2022 12:28"
$ws.Cells.Item(10, 6).Value = "https://www.plus2net.com"
